# Refined Evaluation to be more exact
#
# The predicate-evaluation metrics are split into two finer-grained pairs:
#   - "Correct Extracted Predicates with Parents"  -> renamed to "Detected Predicates Doc Parent"
#   - "Correct Extracted Predicates with Related"   -> renamed to "Detected Predicates Doc Related"
# and two brand-new columns are inserted right after them:
#   - "Correct Pred Predicates Parents"
#   - "Correct Pred Predicates Related"
#
# Inserting two whole columns at O:P shifts every existing column from O
# onward two places to the right (O->Q, P->R, ... U->W), which is exactly
# what the target layout needs, while leaving all of that shifted data
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new, blank columns at O:P — everything from the old column O
# onward (O..U) slides right to (Q..W).
$ws.Columns("O:P").Insert()

# Rename the two predicate-with-parents/related headers to their new labels.
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# Header text for the two newly inserted columns.
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# Fill in the new per-row metric values for the inserted columns.
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 2

$ws.Range("O3").Value = 4
$ws.Range("P3").Value = 4

$ws.Range("O4").Value = 3
$ws.Range("P4").Value = 3

$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 3

$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
